# ============================================================
# Edit: add "2022-Q3" sheet with fund-holding data, and update
# the "总计" (Total) summary sheet with the new quarter's row.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Update the "总计" (Total) summary sheet -------------
$ws1 = $wb.Worksheets.Item(1)

# Need a new row 7 (A7) that doesn't exist yet; give it the same
# style as the existing index column (copy format from A6).
$ws1.Range("A6").Copy()
$ws1.Range("A7").PasteSpecial(-4122)

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 12
$ws1.Range("D2").Value = 1.52
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 6
$ws1.Range("D3").Value = 0.9399999999999999
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 8
$ws1.Range("D4").Value = 0.93
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 13
$ws1.Range("D5").Value = 2.82
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 8
$ws1.Range("D6").Value = 2.09
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 1
$ws1.Range("D7").Value = 0.3

# ---- 2. Insert the new "2022-Q3" worksheet -------------------
# Duplicate an existing 14-row sheet (so the copy already has 13
# data rows + header once we trim one row) right before "2022-Q2",
# then rename it and drop the surplus row.
$srcSheet = $wb.Worksheets.Item(5)
$beforeSheet = $wb.Worksheets.Item(2)
$srcSheet.Copy($beforeSheet)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2022-Q3"
$ws2.Rows("14:14").Delete()

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "'090001"
$ws2.Range("B2").ClearFormats()
$ws2.Range("C2").Value = "大成价值增长混合"
$ws2.Range("D2").Value = "'12.69"
$ws2.Range("D2").ClearFormats()
$ws2.Range("E2").Value = "'64.34"
$ws2.Range("E2").ClearFormats()
$ws2.Range("F2").Value = "'3.01"
$ws2.Range("F2").ClearFormats()
$ws2.Range("G2").Value = "'0.3820"
$ws2.Range("G2").ClearFormats()
$ws2.Range("H2").Value = 10
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "'001915"
$ws2.Range("B3").ClearFormats()
$ws2.Range("C3").Value = "宝盈医疗健康沪港深股票"
$ws2.Range("D3").Value = "'5.15"
$ws2.Range("D3").ClearFormats()
$ws2.Range("E3").Value = "'91.59"
$ws2.Range("E3").ClearFormats()
$ws2.Range("F3").Value = "'6.80"
$ws2.Range("F3").ClearFormats()
$ws2.Range("G3").Value = "'0.3502"
$ws2.Range("G3").ClearFormats()
$ws2.Range("H3").Value = 2
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "'200006"
$ws2.Range("B4").ClearFormats()
$ws2.Range("C4").Value = "长城消费增值混合"
$ws2.Range("D4").Value = "'5.46"
$ws2.Range("D4").ClearFormats()
$ws2.Range("E4").Value = "'90.90"
$ws2.Range("E4").ClearFormats()
$ws2.Range("F4").Value = "'4.11"
$ws2.Range("F4").ClearFormats()
$ws2.Range("G4").Value = "'0.2244"
$ws2.Range("G4").ClearFormats()
$ws2.Range("H4").Value = 4
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "'090020"
$ws2.Range("B5").ClearFormats()
$ws2.Range("C5").Value = "大成健康产业混合A"
$ws2.Range("D5").Value = "'2.38"
$ws2.Range("D5").ClearFormats()
$ws2.Range("E5").Value = "'92.81"
$ws2.Range("E5").ClearFormats()
$ws2.Range("F5").Value = "'8.87"
$ws2.Range("F5").ClearFormats()
$ws2.Range("G5").Value = "'0.2111"
$ws2.Range("G5").ClearFormats()
$ws2.Range("H5").Value = 1
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "'012045"
$ws2.Range("B6").ClearFormats()
$ws2.Range("C6").Value = "大成医药健康股票A"
$ws2.Range("D6").Value = "'1.95"
$ws2.Range("D6").ClearFormats()
$ws2.Range("E6").Value = "'92.51"
$ws2.Range("E6").ClearFormats()
$ws2.Range("F6").Value = "'9.09"
$ws2.Range("F6").ClearFormats()
$ws2.Range("G6").Value = "'0.1773"
$ws2.Range("G6").ClearFormats()
$ws2.Range("H6").Value = 1
$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "'001365"
$ws2.Range("B7").ClearFormats()
$ws2.Range("C7").Value = "大成正向回报灵活配置混合"
$ws2.Range("D7").Value = "'0.57"
$ws2.Range("D7").ClearFormats()
$ws2.Range("E7").Value = "'92.44"
$ws2.Range("E7").ClearFormats()
$ws2.Range("F7").Value = "'8.64"
$ws2.Range("F7").ClearFormats()
$ws2.Range("G7").Value = "'0.0492"
$ws2.Range("G7").ClearFormats()
$ws2.Range("H7").Value = 2
$ws2.Range("A8").Value = 6
$ws2.Range("B8").Value = "'014121"
$ws2.Range("B8").ClearFormats()
$ws2.Range("C8").Value = "大成品质医疗股票A"
$ws2.Range("D8").Value = "'0.57"
$ws2.Range("D8").ClearFormats()
$ws2.Range("E8").Value = "'91.58"
$ws2.Range("E8").ClearFormats()
$ws2.Range("F8").Value = "'8.08"
$ws2.Range("F8").ClearFormats()
$ws2.Range("G8").Value = "'0.0461"
$ws2.Range("G8").ClearFormats()
$ws2.Range("H8").Value = 2
$ws2.Range("A9").Value = 7
$ws2.Range("B9").Value = "'005293"
$ws2.Range("B9").ClearFormats()
$ws2.Range("C9").Value = "诺德新旺灵活配置混合"
$ws2.Range("D9").Value = "'0.38"
$ws2.Range("D9").ClearFormats()
$ws2.Range("E9").Value = "'93.45"
$ws2.Range("E9").ClearFormats()
$ws2.Range("F9").Value = "'7.38"
$ws2.Range("F9").ClearFormats()
$ws2.Range("G9").Value = "'0.0280"
$ws2.Range("G9").ClearFormats()
$ws2.Range("H9").Value = 6
$ws2.Range("A10").Value = 8
$ws2.Range("B10").Value = "'540007"
$ws2.Range("B10").ClearFormats()
$ws2.Range("C10").Value = "汇丰晋信中小盘股票"
$ws2.Range("D10").Value = "'0.56"
$ws2.Range("D10").ClearFormats()
$ws2.Range("E10").Value = "'92.05"
$ws2.Range("E10").ClearFormats()
$ws2.Range("F10").Value = "'4.21"
$ws2.Range("F10").ClearFormats()
$ws2.Range("G10").Value = "'0.0236"
$ws2.Range("G10").ClearFormats()
$ws2.Range("H10").Value = 3
$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = "'012046"
$ws2.Range("B11").ClearFormats()
$ws2.Range("C11").Value = "大成医药健康股票C"
$ws2.Range("D11").Value = "'0.22"
$ws2.Range("D11").ClearFormats()
$ws2.Range("E11").Value = "'92.51"
$ws2.Range("E11").ClearFormats()
$ws2.Range("F11").Value = "'9.09"
$ws2.Range("F11").ClearFormats()
$ws2.Range("G11").Value = "'0.0200"
$ws2.Range("G11").ClearFormats()
$ws2.Range("H11").Value = 1
$ws2.Range("A12").Value = 10
$ws2.Range("B12").Value = "'014122"
$ws2.Range("B12").ClearFormats()
$ws2.Range("C12").Value = "大成品质医疗股票C"
$ws2.Range("D12").Value = "'0.08"
$ws2.Range("D12").ClearFormats()
$ws2.Range("E12").Value = "'91.58"
$ws2.Range("E12").ClearFormats()
$ws2.Range("F12").Value = "'8.08"
$ws2.Range("F12").ClearFormats()
$ws2.Range("G12").Value = "'0.0065"
$ws2.Range("G12").ClearFormats()
$ws2.Range("H12").Value = 2
$ws2.Range("A13").Value = 11
$ws2.Range("B13").Value = "'016060"
$ws2.Range("B13").ClearFormats()
$ws2.Range("C13").Value = "大成健康产业混合C"
$ws2.Range("D13").Value = "'0.00"
$ws2.Range("D13").ClearFormats()
$ws2.Range("E13").Value = "'92.81"
$ws2.Range("E13").ClearFormats()
$ws2.Range("F13").Value = "'8.87"
$ws2.Range("F13").ClearFormats()
$ws2.Range("G13").Value = 0
$ws2.Range("H13").Value = 1
